$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the d5ea74b1-...md file; its status
# moves from "In Translation" to "Ready for handoff" now that the
# handoff package has been generated for both locales.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: same row/file, plus the new handoff timestamp.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "2016-03-08 16:12:56"

# --- de-de sheet: same row/file, plus its own handoff timestamp.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "2016-03-08 16:13:02"
